# Roster update: refresh player/position/team data for rows 9-19
# (P.J. Washington's position gains C eligibility; several players are
# swapped out/in; the table shrinks by one row as Paul George's old row
# is removed after everything shifts up).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# P.J. Washington now also eligible at C
$ws.Range("B9").Value = "PF,C"

# Rows 10-18 get the refreshed roster (name, position, team)
$ws.Range("A10").Value = "Zion Williamson"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "New Orleans Pelicans"

$ws.Range("A11").Value = "Moussa Diabate"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Charlotte Hornets"

$ws.Range("A12").Value = "Ivica Zubac"
$ws.Range("B12").Value = "C"
$ws.Range("C12").Value = "LA Clippers"

$ws.Range("A13").Value = "Goga Bitadze"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "Orlando Magic"

$ws.Range("A14").Value = "Anfernee Simons"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Portland Trail Blazers"

$ws.Range("A15").Value = "Keyonte George"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "Utah Jazz"

$ws.Range("A16").Value = "Giannis Antetokounmpo"
$ws.Range("B16").Value = "PF,C"
$ws.Range("C16").Value = "Milwaukee Bucks"

$ws.Range("A17").Value = "Jonathan Kuminga"
$ws.Range("B17").Value = "SF,PF"
$ws.Range("C17").Value = "Golden State Warriors"

$ws.Range("A18").Value = "Paul George"
$ws.Range("B18").Value = "SG,SF,PF"
$ws.Range("C18").Value = "Philadelphia 76ers"

# The old row 19 (previous Paul George row) is no longer needed
$ws.Rows(19).Delete()
